$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '21.741.67'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.540.45'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.29'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3934'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3205'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.97'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07206'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.077'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9995'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.770'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.58'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.646'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001101'
$ws.Range('E16').Value = '  -2.50%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.462.41'
$ws.Range('E17').Value = '  -5.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06621'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '84.45'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9992'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.167'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.63'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.89'
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.368'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '21.720.89'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.402'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.83'
$ws.Range('E27').Value = '  +2.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.55'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.865'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.711.82'
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.98'
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.139'
$ws.Range('E32').Value = '  +7.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9777'
$ws.Range('E33').Value = '  -7.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08151'
$ws.Range('E34').Value = '  -2.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.611'
$ws.Range('E35').Value = '  -6.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.235'
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02249'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('B38').Value = 'WEMIXTOKEN'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.494'
$ws.Range('E38').Value = '  -6.71%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06022'
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.35'
$ws.Range('E40').Value = '  +6.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2056'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.184'
$ws.Range('E42').Value = '  -2.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9993'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5846'
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.22'
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.734'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5607'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.901'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.169'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '116.44'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06744'
$ws.Range('E51').Value = '  -1.27%  '
